$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as text,
# matching the source data (prices formatted as plain text strings),
# then restore the default "Normal" style so no stray number format
# is left on the cell.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '71.025.44'
$ws.Range('E2').Value = '  +6.33%  '

$ws.Range('D3').Value = '3.658.71'
$ws.Range('E3').Value = '  +6.43%  '

Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.04%  '

Set-TextValue $ws.Range('D5') '596.47'
$ws.Range('E5').Value = '  +2.88%  '

Set-TextValue $ws.Range('D6') '194.69'
$ws.Range('E6').Value = '  +3.44%  '

$ws.Range('D8').Value = '3.651.32'
$ws.Range('E8').Value = '  +6.40%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('E10').Value = '  +8.23%  '

$ws.Range('E11').Value = '  +5.05%  '

Set-TextValue $ws.Range('D12') '58.51'
$ws.Range('E12').Value = '  +3.15%  '

Set-TextValue $ws.Range('D13') '0.0000295'
$ws.Range('E13').Value = '  +6.76%  '

Set-TextValue $ws.Range('D14') '9.99'
$ws.Range('E14').Value = '  +6.20%  '

$ws.Range('D15').Value = '4.243.72'
$ws.Range('E15').Value = '  +6.28%  '

Set-TextValue $ws.Range('D16') '20.09'
$ws.Range('E16').Value = '  +7.41%  '

$ws.Range('D17').Value = '3.657.70'
$ws.Range('E17').Value = '  +6.20%  '

$ws.Range('D18').Value = '71.023.58'
$ws.Range('E18').Value = '  +6.28%  '

$ws.Range('E19').Value = '  +6.66%  '

$ws.Range('E20').Value = '  +2.20%  '

$ws.Range('E21').Value = '  +4.51%  '

Set-TextValue $ws.Range('D22') '489.11'
$ws.Range('E22').Value = '  +1.64%  '

Set-TextValue $ws.Range('D23') '19.09'
$ws.Range('E23').Value = '  +13.41%  '

Set-TextValue $ws.Range('D24') '5.29'
$ws.Range('E24').Value = '  -0.40%  '

Set-TextValue $ws.Range('D25') '4.51'
$ws.Range('E25').Value = '  +4.28%  '

Set-TextValue $ws.Range('D26') '91.58'
$ws.Range('E26').Value = '  +2.73%  '

$ws.Range('E27').Value = '  +7.12%  '

Set-TextValue $ws.Range('D28') '11.46'
$ws.Range('E28').Value = '  +4.62%  '

Set-TextValue $ws.Range('D29') '9.62'
$ws.Range('E29').Value = '  +6.73%  '

$ws.Range('E30').Value = '  +5.69%  '

$ws.Range('E31').Value = '  +5.86%  '

$ws.Range('E32').Value = '  +10.01%  '

Set-TextValue $ws.Range('D33') '630.04'
$ws.Range('E33').Value = '  +5.67%  '

$ws.Range('E34').Value = '  +4.72%  '

Set-TextValue $ws.Range('D35') '67.04'
$ws.Range('E35').Value = '  +4.65%  '

Set-TextValue $ws.Range('D36') '40.33'
$ws.Range('E36').Value = '  +9.90%  '

$ws.Range('D37').Value = '0.0₃0833'
$ws.Range('E37').Value = '  +10.71%  '

Set-TextValue $ws.Range('D38') '0.414'
$ws.Range('E38').Value = '  +7.63%  '

$ws.Range('E39').Value = '  +0.40%  '

$ws.Range('E40').Value = '  +0.15%  '

$ws.Range('E41').Value = '  +2.52%  '

$ws.Range('D42').Value = '3.319.36'
$ws.Range('E42').Value = '  +2.69%  '

$ws.Range('E43').Value = '  +9.81%  '

Set-TextValue $ws.Range('D44') '2.83'
$ws.Range('E44').Value = '  +12.43%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D45') '3.06'
$ws.Range('E45').Value = '  +9.08%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D46') '0.0457'
$ws.Range('E46').Value = '  +6.40%  '

Set-TextValue $ws.Range('D47') '9.50'
$ws.Range('E47').Value = '  +10.47%  '

Set-TextValue $ws.Range('D48') '0.140'
$ws.Range('E48').Value = '  +4.03%  '

$ws.Range('E49').Value = '  +3.25%  '

Set-TextValue $ws.Range('D50') '3.28'
$ws.Range('E50').Value = '  -2.63%  '

Set-TextValue $ws.Range('D51') '0.999'
$ws.Range('E51').Value = '  -0.04%  '
